# Apply "final N-policy linking plots pt 2" edit:
# Add a third iteration block (Iteration_3) in columns K:M, mirroring the
# existing Iteration_1 (E:G) / Iteration_2 (H:J) blocks, and update a few
# existing values in columns E, F, G, H, I, J to reflect the new
# (re-)optimized iteration results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Merge the new header cell exactly like the other iteration headers.
# Do this before touching formatting/values so the existing cell styles
# used for the copy/paste below are not disturbed.
$ws.Range("K1:M1").Merge()

# Force the year labels to be stored as text (matching "2030"/"2040"/"2050"
# as used for the other iteration blocks) instead of being auto-converted
# to numbers.
$ws.Range("K2:M2").NumberFormat = "@"
$ws.Range("K2").Value = "2030"
$ws.Range("L2").Value = "2040"
$ws.Range("M2").Value = "2050"

# --- 2. Bring over header formatting (style only) for the new columns ---
# Row 1 header band (merged "Iteration_2" cell H1:J1 -> new K1:M1)
$ws.Range("H1:J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)

# Row 2 header band (year labels H2:J2 -> new K2:M2)
$ws.Range("H2:J2").Copy()
$ws.Range("K2").PasteSpecial(-4122)

# --- 3. Header text for the new Iteration_3 block ---
$ws.Range("K1").Value = "Iteration_3"

# --- 4. Updated values in existing Iteration_1 / Iteration_2 columns ---
$ws.Range("E4").Value = 1184000.000000019
$ws.Range("H4").Value = 1183999.999999995

$e6 = [double]"1.80409352167688e-08"
$ws.Range("E6").Value = $e6
$ws.Range("F6").Value = 1181739.221044932
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 1181734.682806534

$e7 = [double]"1.979831757394561e-10"
$ws.Range("E7").Value = $e7
$ws.Range("G7").Value = 1183984.965423797
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 1183896.153265644

# --- 5. New Iteration_3 data values (columns K/L/M, rows 4-15) ---
$ws.Range("K4").Value = 1184000.00000003
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0

$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0

$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 1181739.385345381
$ws.Range("M6").Value = 0

$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 1183897.260867802

$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 0

$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 0

$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 0

$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 0

$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 0

$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0

$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 0

$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0

Write-Output "edit applied"
